# DATITORNEO.xlsx edit script
# Summary of changes (per commit message: "PWA: manifest e service worker,
# fix numeri interi, aggiorno Excel e PNG loghi"):
#  - Team "Officina Ronconi" (logo ronconi.png) renamed to "PM Sport"
#    (logo piemme.png) everywhere it is referenced.
#  - "Classifiche" sheet: fixed team-name typo "Commerciali Ferramenta" ->
#    "Commerciale Ferramenta".
#  - "Classifiche" sheet: fixed integer position numbers for Girone B rows
#    (they continued 5..9 instead of restarting at 1..5).
#  - "Classifiche" sheet: removed the stray column L ("√" marks) and left a
#    stray "s" value in D16 (checkmark character is replaced by plain "s").
#  - View/selection state: "Rose" becomes the active/selected sheet, with
#    updated zoom + selected cell; "Classifiche" and "Calendario Partite"
#    get updated zoom/selection too.

$wb = $excel.ActiveWorkbook

$wsProssime    = $wb.Worksheets.Item("Prossime Partite")
$wsRose        = $wb.Worksheets.Item("Rose")
$wsClassifiche = $wb.Worksheets.Item("Classifiche")
$wsCalendario  = $wb.Worksheets.Item("Calendario Partite")

# ---------------------------------------------------------------------
# 1) Rename team "Officina Ronconi" -> "PM Sport" and its logo filename
#    "ronconi.png" -> "piemme.png" everywhere they are used.
# ---------------------------------------------------------------------

# Sheet "Rose"
$wsRose.Range("A4").Value = "PM Sport"

# Sheet "Classifiche"
$wsClassifiche.Range("B4").Value = "PM Sport"
$wsClassifiche.Range("K4").Value = "piemme.png"

# Sheet "Calendario Partite" (3 fixtures involve this team)
$wsCalendario.Range("D3").Value = "PM Sport"
$wsCalendario.Range("E3").Value = "piemme.png"

$wsCalendario.Range("D8").Value = "PM Sport"
$wsCalendario.Range("E8").Value = "piemme.png"

$wsCalendario.Range("F16").Value = "PM Sport"
$wsCalendario.Range("G16").Value = "piemme.png"

# ---------------------------------------------------------------------
# 2) Fix team-name typo on "Classifiche": "Commerciali Ferramenta" ->
#    "Commerciale Ferramenta" (matches spelling already used elsewhere,
#    e.g. "Calendario Partite").
# ---------------------------------------------------------------------
$wsClassifiche.Range("B2").Value = "Commerciale Ferramenta"

# ---------------------------------------------------------------------
# 3) Fix integer position numbers for Girone B rows on "Classifiche":
#    they must restart at 1 instead of continuing from Girone A's count.
# ---------------------------------------------------------------------
$wsClassifiche.Range("A6").Value = 1
$wsClassifiche.Range("A7").Value = 2
$wsClassifiche.Range("A8").Value = 3
$wsClassifiche.Range("A9").Value = 4
$wsClassifiche.Range("A10").Value = 5

# ---------------------------------------------------------------------
# 4) Remove stray column L ("√" marks) on "Classifiche" and leave the
#    lone "s" value behind in D16.
# ---------------------------------------------------------------------
$wsClassifiche.Range("L2").ClearContents()
$wsClassifiche.Range("L3").ClearContents()
$wsClassifiche.Range("D16").Value = "s"

# ---------------------------------------------------------------------
# 5) View / selection state updates.
# ---------------------------------------------------------------------

# "Calendario Partite": zoom 132 -> 112, selection D2 -> G2
$wsCalendario.Activate()
$excel.ActiveWindow.Zoom = 112
$wsCalendario.Range("G2").Select()

# "Classifiche": zoom stays 169, selection L4 -> B2 (no longer the tab
# that is selected/active when the workbook is reopened).
$wsClassifiche.Activate()
$excel.ActiveWindow.Zoom = 169
$wsClassifiche.Range("B2").Select()

# "Rose" becomes the active/selected sheet: zoom 118 -> 237, selection
# E15 -> A5. Activated last so it ends up as the workbook's active tab.
$wsRose.Activate()
$excel.ActiveWindow.Zoom = 237
$wsRose.Range("A5").Select()
